$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Row 88: finished frames (Rahmen) log entry ---
$ws.Range("A88").Value2 = 43219
$ws.Range("C88").Value = "Ein Massiver und ein Leichter Ramen. Das Testen kann beginnen!"
$ws.Range("B88").Value = "2018-04-29 1.JPG"

# --- Row 89: two-way communication working ---
$ws.Range("A89").Value2 = 43219
$ws.Range("B89").Value = "2018-04-29 2.AVI"
$ws.Range("C89").Value = "Die Kommunikation hin und zurück läuft! Drücken auf Start -> Biep"

# Reuse the same date format/style as the rest of column A (avoids minting a
# new style entry).
$ws.Range("A87").Copy()
$ws.Range("A88:A89").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Mirror the view/selection state after the new rows were added
$ws.Range("A90").Select()
